$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 7, pushing the existing data (old rows 7-21)
# down to rows 9-23, to make room for the new weekly report rows.
$ws.Rows("7:8").Insert()

# New row 7: Camote, "1a (cosecha)"
$ws.Cells.Item(7, 1).Value = 1
$ws.Cells.Item(7, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(7, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(7, 4).Value = 44547
$ws.Cells.Item(7, 5).Value = 15
$ws.Cells.Item(7, 6).Value = 100112045
$ws.Cells.Item(7, 7).Value = "Zapallo"
$ws.Cells.Item(7, 8).Value = "Camote"
$ws.Cells.Item(7, 9).Value = "1a (cosecha)"
$ws.Cells.Item(7, 10).Value = 800
$ws.Cells.Item(7, 11).Value = 600
$ws.Cells.Item(7, 12).Value = 650
$ws.Cells.Item(7, 13).Value = 625
$ws.Cells.Item(7, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(7, 15).Value = "Perú"
$ws.Cells.Item(7, 16).Value = 625
$ws.Cells.Item(7, 17).Value = 1
$ws.Cells.Item(7, 18).Value = "Hortaliza"

# New row 8: Camote, "2a nueva(o)"
$ws.Cells.Item(8, 1).Value = 1
$ws.Cells.Item(8, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(8, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(8, 4).Value = 44547
$ws.Cells.Item(8, 5).Value = 15
$ws.Cells.Item(8, 6).Value = 100112045
$ws.Cells.Item(8, 7).Value = "Zapallo"
$ws.Cells.Item(8, 8).Value = "Camote"
$ws.Cells.Item(8, 9).Value = "2a nueva(o)"
$ws.Cells.Item(8, 10).Value = 300
$ws.Cells.Item(8, 11).Value = 500
$ws.Cells.Item(8, 12).Value = 550
$ws.Cells.Item(8, 13).Value = 525
$ws.Cells.Item(8, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(8, 15).Value = "Perú"
$ws.Cells.Item(8, 16).Value = 525
$ws.Cells.Item(8, 17).Value = 1
$ws.Cells.Item(8, 18).Value = "Hortaliza"
